# Simplify the header/column names (sharedStrings text updates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Pending"
$ws.Range("C1").Value = "Positive"
$ws.Range("D1").Value = "Negative"
$ws.Range("E1").Value = "Total_tested"
$ws.Range("F1").Value = "In_quarantine"
$ws.Range("G1").Value = "Released_quarantine"

# Move the active selection to G2 (matches the saved sheetView selection)
$ws.Range("G2").Select()

# Widen column G (target OOXML width 19.64453125 chars) and shrink row 1's
# height (was wrapping to a taller row with the long header text, now fits
# in a shorter row with the short header text).
$ws.Columns.Item(7).ColumnWidth = 18.83
$ws.Rows.Item(1).RowHeight = 28.7
